$p = $ppt.ActivePresentation

# 1) Merge the split "What does MVS from the " / "E-Land " / "toolbox do?" runs
#    inside the SmartArt diagram on slide 5 into a single run with the
#    combined text. Updating the SmartArt node text rewrites both the
#    diagram data (data1.xml) and the cached drawing (drawing1.xml).
$s5 = $p.Slides.Item(5)
$diagramShape = $s5.Shapes.Item(4)
$smartArt = $diagramShape.SmartArt
for ($i = 1; $i -le $smartArt.AllNodes.Count; $i++) {
    $node = $smartArt.AllNodes.Item($i)
    if ($node.TextFrame2.TextRange.Text -eq "What does MVS from the E-Land toolbox do?") {
        $node.TextFrame2.TextRange.Text = "What does MVS from the E-Land toolbox do?"
    }
}

# 2) Fix typo "stundets" -> "students" on slide 3
$s3 = $p.Slides.Item(3)
$textShape = $s3.Shapes.Item(4)
[void]$textShape.TextFrame.TextRange.Replace("stundets", "students")
